{"js": "// 1. Update \"{{ test_date | date(format=\u201d%Y\u201d) }}\" -> \"{{ test_date | date(format=\u201dY\u201d) }}\"\n//    and insert a new paragraph right after it for the new `test_time | time(...)` filter.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet dateFormatParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"{{ test_date | date(format=\") === 0) {\n    dateFormatParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!dateFormatParagraph) {\n  throw new Error(\"Could not find the 'test_date | date(format=...)' paragraph\");\n}\n\n// Narrow the \"%Y\" inside that paragraph's range so we only touch this occurrence.\nconst yearRange = dateFormatParagraph.search(\"%Y\", { matchCase: true });\nyearRange.load(\"items/text\");\nawait context.sync();\nif (yearRange.items.length === 0) {\n  throw new Error(\"Could not find '%Y' inside the date-format paragraph\");\n}\nyearRange.items[0].insertText(\"Y\", \"Replace\");\nawait context.sync();\n\n// Insert the new paragraph for the time filter right after the date-format paragraph.\ndateFormatParagraph.insertParagraph(\n  \"{{ test_time | time(format=\u201dH:m\u201d) }}\",\n  \"After\"\n);\nawait context.sync();\n\n// 2. In the \"datetime2\" paragraph (syntax-highlighted code sample), rename the\n//    \"in_format\" keyword run to \"format\" and drop the whole \"in_format=...out_format=\"\n//    sub-range, keeping only the trailing format string (now also updated).\nconst inFormatRange = body.search(\"in_format\", { matchCase: true });\ninFormatRange.load(\"items/text\");\nawait context.sync();\nif (inFormatRange.items.length === 0) {\n  throw new Error(\"Could not find the 'in_format' run\");\n}\ninFormatRange.items[0].insertText(\"format\", \"Replace\");\nawait context.sync();\n\nconst droppedRange = body.search('\"%H:%M-%Y-%m-%d\", out_format=', {\n  matchCase: true,\n});\ndroppedRange.load(\"items/text\");\nawait context.sync();\nif (droppedRange.items.length === 0) {\n  throw new Error(\"Could not find the text range to drop\");\n}\ndroppedRange.items[0].delete();\nawait context.sync();\n\nconst outFormatValueRange = body.search('\"%H:%M\u201d)', { matchCase: true });\noutFormatValueRange.load(\"items/text\");\nawait context.sync();\nif (outFormatValueRange.items.length === 0) {\n  throw new Error(\"Could not find the '\\\"%H:%M\u201d)' run\");\n}\noutFormatValueRange.items[0].insertText('\"H:m\u201d)', \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Helper characters (avoid any literal curly-quote encoding issues in source).\n$quote = [char]34        # \"  (straight double quote)\n$curly = [char]8221      # \u201d  (right double quotation mark, U+201D)\n\n# ---------------------------------------------------------------------------\n# 1. Locate the \"{{ test_date | date(format=...) }}\" paragraph and change the\n#    strftime-style \"%Y\" token to the new bare \"Y\" token.\n# ---------------------------------------------------------------------------\n$dateFormatParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"{{ test_date | date(format=\")) {\n        $dateFormatParagraph = $p\n        break\n    }\n}\nif ($null -eq $dateFormatParagraph) {\n    throw \"Could not find the 'test_date | date(format=...)' paragraph\"\n}\n\n$yearRange = $dateFormatParagraph.Range\n$yearFound = $yearRange.Find.Execute(\"%Y\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $yearFound) {\n    throw \"Could not find '%Y' inside the date-format paragraph\"\n}\n$yearRange.Text = \"Y\"\n\n# ---------------------------------------------------------------------------\n# 2. Insert a brand-new paragraph right after it for the new `time` filter.\n# ---------------------------------------------------------------------------\n$insertionRange = $dateFormatParagraph.Range\n$insertionRange.InsertParagraphAfter()\n\n$newParagraph = $dateFormatParagraph.Next()\n$newRange = $newParagraph.Range\n$newRange.InsertAfter(\"{{ test_time | time(format=\" + $curly + \"H:m\" + $curly + \") }}\")\n\n# ---------------------------------------------------------------------------\n# 3. In the \"datetime2\" paragraph (syntax-highlighted code sample):\n#      a. rename the \"in_format\" keyword run to \"format\"\n#      b. drop the whole in-format value / \", \" / \"out_format\" / \"=\" runs\n#      c. shorten the remaining format-string run from \"%H:%M\" to \"H:m\"\n# ---------------------------------------------------------------------------\n$datetime2Paragraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"{{ test_datetime2 |\")) {\n        $datetime2Paragraph = $p\n        break\n    }\n}\nif ($null -eq $datetime2Paragraph) {\n    throw \"Could not find the 'test_datetime2 | datetime(...)' paragraph\"\n}\n\n$inFormatRange = $datetime2Paragraph.Range\n$inFormatFound = $inFormatRange.Find.Execute(\"in_format\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $inFormatFound) {\n    throw \"Could not find the 'in_format' run\"\n}\n$inFormatRange.Text = \"format\"\n\n$dropRange = $datetime2Paragraph.Range\n$dropTarget = $quote + \"%H:%M-%Y-%m-%d\" + $quote + \", out_format=\"\n$dropFound = $dropRange.Find.Execute($dropTarget, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $dropFound) {\n    throw \"Could not find the text range to drop\"\n}\n$dropRange.Text = \"\"\n\n$valueRange = $datetime2Paragraph.Range\n$valueFind = $quote + \"%H:%M\" + $curly + \")\"\n$valueFound = $valueRange.Find.Execute($valueFind, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $valueFound) {\n    throw \"Could not find the quoted '%H:%M' format-string run\"\n}\n$valueRange.Text = $quote + \"H:m\" + $curly + \")\"\n"}
